# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the cryptos worksheet.
# D-column cells are forced to Text format before assignment so that
# numeric-looking price strings (e.g. "111.34") are stored as text,
# matching the original "inlineStr" cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.491.69'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.630.77'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.34'
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '325.28'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  -1.43%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.37'
$ws.Range("E10").Value = '  -4.37%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0809'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.35'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.038.03'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.638.78'
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("E17").Value = '  -2.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.409.73'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.67'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '267.08'
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.85'
$ws.Range("E24").Value = '  -4.86%  '
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("E28").Value = '  +2.15%  '
$ws.Range("E29").Value = '  -3.07%  '
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.55'
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.61'
$ws.Range("E32").Value = '  -1.45%  '
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0808'
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.93'
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.95'
$ws.Range("E37").Value = '  +3.77%  '
$ws.Range("E38").Value = '  -2.97%  '
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.64'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.78'
$ws.Range("E41").Value = '  +1.55%  '
$ws.Range("E42").Value = '  -1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.21'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0326'
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.041.47'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +8.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.22'
$ws.Range("E47").Value = '  -3.79%  '
$ws.Range("E48").Value = '  -4.10%  '
$ws.Range("E49").Value = '  -3.71%  '
$ws.Range("E50").Value = '  -3.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.45'
$ws.Range("E51").Value = '  +1.21%  '
